# Auto-generated edit script: updates computed market-profit columns (H-N)
# across multiple Leve sheets per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11 (Leve Item ID 5533)
$ws.Range("H11").Value = 358.17648
$ws.Range("I11").Value = 358.17648
$ws.Range("K11").Value = 358.17648
$ws.Range("M11").Value = -218.17648
# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 1519.619
$ws.Range("I19").Value = 594.875
$ws.Range("J19").Value = 4478.8
$ws.Range("K19").Value = 594.875
$ws.Range("L19").Value = 4478.8
$ws.Range("M19").Value = -419.875
$ws.Range("N19").Value = -4828.8
# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 2200
$ws.Range("I41").Value = 2900
$ws.Range("K41").Value = 2900
$ws.Range("M41").Value = -2460
# Row 51 (Leve Item ID 5486)
$ws.Range("H51").Value = 13892973
$ws.Range("I51").Value = 4750
$ws.Range("J51").Value = 20837084
$ws.Range("K51").Value = 4750
$ws.Range("L51").Value = 20837084
$ws.Range("N51").Value = -20838052
$ws.Range("M51").Value = -4266
# Row 52 (Leve Item ID 4567)
$ws.Range("H52").Value = 799.5
$ws.Range("I52").Value = 799.5
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 2398.5
$ws.Range("L52").Value = 0
$ws.Range("M52").Value = -2238.5
$ws.Range("N52").ClearContents()
# Row 107 (Leve Item ID 27766)
$ws.Range("H107").Value = 1056.4166
$ws.Range("I107").Value = 1056.4166
$ws.Range("K107").Value = 1056.4166
$ws.Range("M107").Value = 863.5834
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 4878.825
$ws.Range("I137").Value = 3740.16
$ws.Range("J137").Value = 6776.6
$ws.Range("K137").Value = 11220.48
$ws.Range("L137").Value = 20329.8
$ws.Range("M137").Value = -8670.48
$ws.Range("N137").Value = -25429.8
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2668.9583
$ws.Range("I138").Value = 1801
$ws.Range("J138").Value = 2958.2778
$ws.Range("K138").Value = 5403
$ws.Range("L138").Value = 8874.8334
$ws.Range("M138").Value = -263
$ws.Range("N138").Value = -19154.8334

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Range("H2").Value = 3422.0938
$ws.Range("I2").Value = 3673.0417
$ws.Range("K2").Value = 3673.0417
$ws.Range("M2").Value = -3560.0417
# Row 63 (Leve Item ID 12528)
$ws.Range("H63").Value = 1880.6666
$ws.Range("J63").Value = 1922
$ws.Range("L63").Value = 1922
$ws.Range("N63").Value = -3294
# Row 66 (Leve Item ID 12528)
$ws.Range("H66").Value = 1880.6666
$ws.Range("J66").Value = 1922
$ws.Range("L66").Value = 9610
$ws.Range("N66").Value = -16474
# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 1649.6111
$ws.Range("I74").Value = 1411.7693
$ws.Range("J74").Value = 2268
$ws.Range("K74").Value = 1411.7693
$ws.Range("L74").Value = 2268
$ws.Range("M74").Value = -537.7692999999999
$ws.Range("N74").Value = -4016
# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 1649.6111
$ws.Range("I77").Value = 1411.7693
$ws.Range("J77").Value = 2268
$ws.Range("K77").Value = 7058.8465
$ws.Range("L77").Value = 11340
$ws.Range("M77").Value = -2690.8465
$ws.Range("N77").Value = -20076
# Row 102 (Leve Item ID 19945)
$ws.Range("H102").Value = 2272.8462
$ws.Range("I102").Value = 1775.091
$ws.Range("K102").Value = 1775.091
$ws.Range("M102").Value = -153.0909999999999
# Row 116 (Leve Item ID 27713)
$ws.Range("H116").Value = 3422.0938
$ws.Range("I116").Value = 3673.0417
$ws.Range("K116").Value = 3673.0417
$ws.Range("M116").Value = -1379.0417

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Range("H3").Value = 3422.0938
$ws.Range("I3").Value = 3673.0417
$ws.Range("K3").Value = 3673.0417
$ws.Range("M3").Value = -3559.0417
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 3493.8572
$ws.Range("I99").Value = 3678
$ws.Range("K99").Value = 3678
$ws.Range("M99").Value = -2180

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2953.5938
$ws.Range("I31").Value = 2133.8
$ws.Range("J31").Value = 3326.2273
$ws.Range("K31").Value = 2133.8
$ws.Range("L31").Value = 3326.2273
$ws.Range("M31").Value = -1838.8
$ws.Range("N31").Value = -3916.2273
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2953.5938
$ws.Range("I34").Value = 2133.8
$ws.Range("J34").Value = 3326.2273
$ws.Range("K34").Value = 2133.8
$ws.Range("L34").Value = 3326.2273
$ws.Range("M34").Value = -1931.8
$ws.Range("N34").Value = -3730.2273
# Row 39 (Leve Item ID 1915)
$ws.Range("H39").Value = 3350
$ws.Range("I39").Value = 3350
$ws.Range("K39").Value = 3350
$ws.Range("M39").Value = -2959
# Row 49 (Leve Item ID 1915)
$ws.Range("H49").Value = 3350
$ws.Range("I49").Value = 3350
$ws.Range("K49").Value = 3350
$ws.Range("M49").Value = -3168
# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 24623.75
$ws.Range("J68").Value = 24998.572
$ws.Range("L68").Value = 24998.572
$ws.Range("N68").Value = -26496.572
# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 24623.75
$ws.Range("J71").Value = 24998.572
$ws.Range("L71").Value = 74995.716
$ws.Range("N71").Value = -82483.716
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 19196.46
$ws.Range("I99").Value = 25126
$ws.Range("K99").Value = 25126
$ws.Range("M99").Value = -23628
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 19196.46
$ws.Range("I126").Value = 25126
$ws.Range("K126").Value = 75378
$ws.Range("M126").Value = -72908
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 3128.6453
$ws.Range("I134").Value = 2985.4
$ws.Range("J134").Value = 3725.5
$ws.Range("K134").Value = 8956.200000000001
$ws.Range("L134").Value = 11176.5
$ws.Range("M134").Value = -6421.200000000001
$ws.Range("N134").Value = -16246.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 80 (Leve Item ID 12890)
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83 (Leve Item ID 12890)
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# Row 98 (Leve Item ID 19843)
$ws.Range("H98").Value = 728
$ws.Range("J98").Value = 978
$ws.Range("L98").Value = 2934
$ws.Range("N98").Value = -5930
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 3875.6875
$ws.Range("I140").Value = 3617.9
$ws.Range("J140").Value = 4305.3335
$ws.Range("K140").Value = 10853.7
$ws.Range("L140").Value = 12916.0005
$ws.Range("M140").Value = -5673.700000000001
$ws.Range("N140").Value = -23276.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 4254.273
$ws.Range("J80").Value = 5600
$ws.Range("L80").Value = 5600
$ws.Range("N80").Value = -7596
# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 4254.273
$ws.Range("J83").Value = 5600
$ws.Range("L83").Value = 28000
$ws.Range("N83").Value = -37984
# Row 102 (Leve Item ID 36169)
$ws.Range("H102").Value = 1884.1666
$ws.Range("I102").Value = 2061
$ws.Range("K102").Value = 2061
$ws.Range("M102").Value = -439
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 2928.2273
$ws.Range("I113").Value = 2308.7778
$ws.Range("K113").Value = 2308.7778
$ws.Range("M113").Value = -138.7777999999998
# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2322
$ws.Range("I126").Value = 2362.25
$ws.Range("K126").Value = 7086.75
$ws.Range("M126").Value = -4616.75
# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 3374.8928
$ws.Range("I132").Value = 3537.5908
$ws.Range("J132").Value = 2778.3333
$ws.Range("K132").Value = 10612.7724
$ws.Range("L132").Value = 8334.999899999999
$ws.Range("M132").Value = -8082.7724
$ws.Range("N132").Value = -13394.9999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 3235.1614
$ws.Range("I136").Value = 2539.3635
$ws.Range("K136").Value = 7618.0905
$ws.Range("M136").Value = -5068.0905

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 29413138
$ws.Range("J107").Value = 62501200
$ws.Range("L107").Value = 187503600
$ws.Range("N107").Value = -187507440
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 722.1875
$ws.Range("I136").Value = 425.2143
$ws.Range("K136").Value = 1275.6429
$ws.Range("M136").Value = 1274.3571
